$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.419687032699585
$ws.Range("B1").Value = 1.898656964302063
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.922520875930786
$ws.Range("E1").Value = 0.7132148146629333
